# The workbook's first sheet ("Sheet1") originally duplicated the numeric
# "id" column (A) alongside the shared-string "desc" column (B). The id
# column is being dropped from this sheet, leaving only column B populated
# for rows 2 through 42 (row 1 is the header and is left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2:A42").Clear()
